# Rename three header cells in row 1 on "Sheet 1":
#   U1: "Volume (nL)\nDMSO normalization"  -> "DMSO"
#   V1: "Volume (nL)\na+Tw normalization"  -> "Tween"
#   X1: "DMSO %"                           -> "DMSO_pct"
# (W1 "Total well volume (nL)" is untouched; the shared-strings table is
# rebuilt on save, which is why its <v> index shifts even though the text
# doesn't change.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("U1").Value = "DMSO"
$ws.Range("V1").Value = "Tween"
$ws.Range("X1").Value = "DMSO_pct"

# Move the selection/active cell the way the author left it after editing.
$ws.Range("U7").Select()
